$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K ("intervention_type") mirrors the header style used by the other header cells (copy format from J1)
$ws.Range("K1").Value = "intervention_type"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# Populate intervention_type per clinical trial row (rows with no source value keep column K blank)
$ws.Cells.Item(2, 11).Value = "DEVICE"
$ws.Cells.Item(3, 11).Value = "DEVICE"
$ws.Cells.Item(4, 11).Value = "DEVICE"
$ws.Cells.Item(5, 11).Value = "DEVICE"
$ws.Cells.Item(7, 11).Value = "DRUG"
$ws.Cells.Item(8, 11).Value = "DRUG"
$ws.Cells.Item(9, 11).Value = "PROCEDURE"
$ws.Cells.Item(10, 11).Value = "DRUG"
$ws.Cells.Item(11, 11).Value = "DEVICE"
$ws.Cells.Item(12, 11).Value = "DRUG"
$ws.Cells.Item(13, 11).Value = "DRUG"
$ws.Cells.Item(14, 11).Value = "DRUG"
$ws.Cells.Item(15, 11).Value = "DRUG"
$ws.Cells.Item(16, 11).Value = "DEVICE"
$ws.Cells.Item(17, 11).Value = "DRUG"
$ws.Cells.Item(18, 11).Value = "OTHER"
$ws.Cells.Item(19, 11).Value = "PROCEDURE"
$ws.Cells.Item(20, 11).Value = "DEVICE"
$ws.Cells.Item(21, 11).Value = "DRUG"
$ws.Cells.Item(22, 11).Value = "PROCEDURE"
$ws.Cells.Item(23, 11).Value = "BIOLOGICAL"
$ws.Cells.Item(24, 11).Value = "DEVICE"
$ws.Cells.Item(26, 11).Value = "OTHER"
$ws.Cells.Item(27, 11).Value = "BIOLOGICAL"
$ws.Cells.Item(28, 11).Value = "BIOLOGICAL"
$ws.Cells.Item(29, 11).Value = "PROCEDURE"
$ws.Cells.Item(30, 11).Value = "DRUG"
$ws.Cells.Item(31, 11).Value = "DEVICE"
$ws.Cells.Item(32, 11).Value = "DRUG"
$ws.Cells.Item(33, 11).Value = "DRUG"
$ws.Cells.Item(34, 11).Value = "PROCEDURE"
$ws.Cells.Item(35, 11).Value = "DEVICE"
$ws.Cells.Item(36, 11).Value = "PROCEDURE"
$ws.Cells.Item(37, 11).Value = "DEVICE"
$ws.Cells.Item(38, 11).Value = "PROCEDURE"
$ws.Cells.Item(39, 11).Value = "OTHER"
$ws.Cells.Item(40, 11).Value = "DRUG"
$ws.Cells.Item(41, 11).Value = "DRUG"
$ws.Cells.Item(42, 11).Value = "OTHER"
$ws.Cells.Item(43, 11).Value = "PROCEDURE"
$ws.Cells.Item(44, 11).Value = "DRUG"
$ws.Cells.Item(45, 11).Value = "OTHER"
$ws.Cells.Item(46, 11).Value = "DRUG"
$ws.Cells.Item(47, 11).Value = "DEVICE"
$ws.Cells.Item(48, 11).Value = "OTHER"
$ws.Cells.Item(49, 11).Value = "DEVICE"
$ws.Cells.Item(50, 11).Value = "BEHAVIORAL"
$ws.Cells.Item(51, 11).Value = "DEVICE"
$ws.Cells.Item(52, 11).Value = "PROCEDURE"
$ws.Cells.Item(53, 11).Value = "BIOLOGICAL"
$ws.Cells.Item(54, 11).Value = "OTHER"
$ws.Cells.Item(55, 11).Value = "OTHER"
$ws.Cells.Item(56, 11).Value = "DEVICE"
$ws.Cells.Item(57, 11).Value = "OTHER"
$ws.Cells.Item(58, 11).Value = "DEVICE"
$ws.Cells.Item(59, 11).Value = "OTHER"
$ws.Cells.Item(60, 11).Value = "DEVICE"
$ws.Cells.Item(61, 11).Value = "DEVICE"
$ws.Cells.Item(62, 11).Value = "BEHAVIORAL"
$ws.Cells.Item(63, 11).Value = "OTHER"
$ws.Cells.Item(64, 11).Value = "PROCEDURE"
$ws.Cells.Item(65, 11).Value = "DEVICE"
$ws.Cells.Item(66, 11).Value = "OTHER"
$ws.Cells.Item(67, 11).Value = "PROCEDURE"
$ws.Cells.Item(68, 11).Value = "OTHER"
$ws.Cells.Item(69, 11).Value = "DEVICE"
$ws.Cells.Item(70, 11).Value = "PROCEDURE"
$ws.Cells.Item(71, 11).Value = "DIAGNOSTIC_TEST"
$ws.Cells.Item(72, 11).Value = "OTHER"
$ws.Cells.Item(73, 11).Value = "DEVICE"
$ws.Cells.Item(74, 11).Value = "COMBINATION_PRODUCT"
$ws.Cells.Item(75, 11).Value = "OTHER"
$ws.Cells.Item(76, 11).Value = "BEHAVIORAL"
$ws.Cells.Item(77, 11).Value = "PROCEDURE"
$ws.Cells.Item(78, 11).Value = "OTHER"
$ws.Cells.Item(79, 11).Value = "DEVICE"
$ws.Cells.Item(80, 11).Value = "PROCEDURE"
$ws.Cells.Item(81, 11).Value = "OTHER"
$ws.Cells.Item(82, 11).Value = "OTHER"
$ws.Cells.Item(83, 11).Value = "DEVICE"
$ws.Cells.Item(84, 11).Value = "OTHER"
$ws.Cells.Item(85, 11).Value = "DIAGNOSTIC_TEST"
$ws.Cells.Item(86, 11).Value = "DEVICE"
$ws.Cells.Item(87, 11).Value = "DIAGNOSTIC_TEST"
$ws.Cells.Item(88, 11).Value = "RADIATION"
$ws.Cells.Item(89, 11).Value = "DEVICE"
$ws.Cells.Item(90, 11).Value = "OTHER"
$ws.Cells.Item(91, 11).Value = "OTHER"
$ws.Cells.Item(92, 11).Value = "DIAGNOSTIC_TEST"
$ws.Cells.Item(93, 11).Value = "DEVICE"
$ws.Cells.Item(94, 11).Value = "DEVICE"
$ws.Cells.Item(95, 11).Value = "DEVICE"
